# Update the station load data (columns B:I, rows 2-18) with the new
# measured values as reflected in the updated source data export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(19189,5219,8740,0,5529,0,15018,19463),
  @(18379,4985,8351,0,5430,0,14349,18477),
  @(17470,4704,8155,0,5207,0,13877,17820),
  @(17438,4638,8102,0,5178,0,13798,17620),
  @(18041,4882,8235,0,5226,0,14143,17900),
  @(18811,5371,8886,0,5719,0,14910,19136),
  @(18689,4586,9006,0,7021,0,14958,18598),
  @(21111,4741,10697,0,10465,0,17456,20215),
  @(27401,5613,13587,0,12996,0,19031,23924),
  @(29301,5688,17571,0,14723,0,20103,24871),
  @(30596,5902,18831,0,14864,0,20562,25402),
  @(30033,6037,18926,0,14319,0,20138,25728),
  @(30101,5877,19038,0,15795,0,21958,25669),
  @(31220,6127,19084,0,15071,0,21170,26296),
  @(32197,6332,19180,0,15521,0,21620,26956),
  @(29317,5863,18637,0,14787,0,20501,25861),
  @(0,5726,0,0,0,0,0,9830)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
